# Add a new "time_taken" column (F) to the worksheet with a header and
# per-row timestamp metadata, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:39:38.943298",
    "2021-10-05 13:39:38.943309",
    "2021-10-05 13:39:38.943313",
    "2021-10-05 13:39:38.943316",
    "2021-10-05 13:39:38.943320",
    "2021-10-05 13:39:38.943323",
    "2021-10-05 13:39:38.943326",
    "2021-10-05 13:39:38.943329",
    "2021-10-05 13:39:38.943332",
    "2021-10-05 13:39:38.943335",
    "2021-10-05 13:39:38.943339",
    "2021-10-05 13:39:38.943342",
    "2021-10-05 13:39:38.943345",
    "2021-10-05 13:39:38.943348",
    "2021-10-05 13:39:38.943351",
    "2021-10-05 13:39:38.943353",
    "2021-10-05 13:39:38.943357",
    "2021-10-05 13:39:38.943360",
    "2021-10-05 13:39:38.943363",
    "2021-10-05 13:39:38.943366",
    "2021-10-05 13:39:38.943369",
    "2021-10-05 13:39:38.943372",
    "2021-10-05 13:39:38.943375",
    "2021-10-05 13:39:38.943378",
    "2021-10-05 13:39:38.943381",
    "2021-10-05 13:39:38.943384"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
